# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update case numbers for several countries (rows sorted desc. by Col B)
# - Costa Rica / Guinea Ecuatorial swap places (Guinea Ecuatorial's updated
#   totals now tie with Costa Rica's, so it moves above Costa Rica in the
#   country ranking)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 22:05"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 1613950
$ws.Cells.Item(4, 3).Value = 21227
$ws.Cells.Item(4, 4).Value = 374588
$ws.Cells.Item(4, 5).Value = 1143344
$ws.Cells.Item(4, 7).Value = 1082
$ws.Cells.Item(4, 8).Value = 96018

# --- Alemania (row 11) ---
$ws.Cells.Item(11, 2).Value = 178918
$ws.Cells.Item(11, 3).Value = 387
$ws.Cells.Item(11, 5).Value = 12636
$ws.Cells.Item(11, 7).Value = 12
$ws.Cells.Item(11, 8).Value = 8282

# --- Peru (row 15) ---
$ws.Cells.Item(15, 2).Value = 108769
$ws.Cells.Item(15, 3).Value = 4749
$ws.Cells.Item(15, 4).Value = 43587
$ws.Cells.Item(15, 5).Value = 62034
$ws.Cells.Item(15, 7).Value = 124
$ws.Cells.Item(15, 8).Value = 3148

# --- Canada (row 17) ---
$ws.Cells.Item(17, 2).Value = 81279
$ws.Cells.Item(17, 3).Value = 1137
$ws.Cells.Item(17, 4).Value = 41623
$ws.Cells.Item(17, 5).Value = 33511

# --- Barein (row 55) ---
$ws.Cells.Item(55, 2).Value = 8188
$ws.Cells.Item(55, 3).Value = 300
$ws.Cells.Item(55, 4).Value = 3873
$ws.Cells.Item(55, 5).Value = 4303

# --- Costa Rica / Guinea Ecuatorial (rows 113 & 114) ---
# Guinea Ecuatorial gets fresh totals and takes row 113 (now tied with
# Costa Rica's 903 total cases, alphabetically/ranking-wise ahead of it);
# Costa Rica (unchanged data) drops to row 114.
$ws.Cells.Item(113, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(113, 2).Value = 903
$ws.Cells.Item(113, 3).Value = 13
$ws.Cells.Item(113, 4).Value = 22
$ws.Cells.Item(113, 5).Value = 871
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = 10

$ws.Cells.Item(114, 1).Value = "Costa Rica"
$ws.Cells.Item(114, 2).Value = 903
$ws.Cells.Item(114, 3).Value = 6
$ws.Cells.Item(114, 4).Value = 592
$ws.Cells.Item(114, 5).Value = 301
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 10

# --- Etiopia (row 139) ---
$ws.Cells.Item(139, 2).Value = 399
$ws.Cells.Item(139, 3).Value = 10
$ws.Cells.Item(139, 5).Value = 271
